# Initiating create engagement button
# Adds a new TODO item row ("Loguer avec le matricule au lieu de l'email")
# right after the first entry on the TODO sheet, pushing the existing
# rows down by one, and highlights the first row the same way row 2 is
# already highlighted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO")

# Highlight A1/B1 with the same "done" fill used on row 2 (A2/B2).
$highlightColor = $ws.Range("A2").Interior.Color()
$ws.Range("A1").Interior.Color = $highlightColor
$ws.Range("B1").Interior.Color = $highlightColor

# Insert a new row before the current row 3, shifting rows 3-20 down to 4-21.
# Excel copies formatting from the row above (row 2), which is the style
# we want for the new A3/B3 cells.
$ws.Rows.Item(3).Insert()

# Populate the new TODO entry.
$ws.Cells.Item(3, 2).Value = "Loguer avec le matricule au lieu de l'email"

# Update the selected cell to reflect where the user was working.
$ws.Range("B4").Select() | Out-Null
